# Rename the sheet (B08L7BD4C6 -> B08MPM2BB2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "B08MPM2BB2"

# Replace the product/keyword list in column A (A1:A100) with the new values
$values = @(
    "black top",
    "black bra",
    "wireless bra",
    "sport bra",
    "pink top",
    "yoga bra",
    "yoga top",
    "low back bra",
    "black s",
    "longline bra",
    "padded bra",
    "green bra",
    "strappy bra",
    "bra top",
    "pink bra",
    "pink bralette",
    "workout bra",
    "strappy bralette",
    "black bra top",
    "longline bralette",
    "black cross",
    "bralette top",
    "cute top",
    "criss cross bra",
    "black l",
    "black m",
    "wirefree bra",
    "criss cross top",
    "long line bra",
    "cute bra",
    "criss cross",
    "pink bra top",
    "bralette bra",
    "impact sport",
    "long s",
    "top s",
    "black cup",
    "wireless bralette",
    "women’s longline sports bra wirefree padded medium support yoga bras gym running workout tank tops",
    "cross bra",
    "green l",
    "cup with",
    "sport gym",
    "long bra",
    "yoga gym",
    "low back bralette",
    "longline yoga bra",
    "criss cross bralette",
    "sport elastic",
    "cute back",
    "pink cross",
    "green s",
    "black apparel",
    "black yoga top",
    "green m",
    "yoga bra strappy",
    "cross back bra",
    "cute pink",
    "yoga workout",
    "sport workout",
    "cute bralette",
    "low back top",
    "criss cross back top",
    "strappy top",
    "criss cross front",
    "strappy sport bra",
    "cross front bra",
    "bra elastic",
    "long line",
    "sport top",
    "strappy back bralette",
    "green sport bra",
    "bra cup",
    "wireless sport",
    "top m",
    "strappy back top",
    "gym apparel",
    "gym back",
    "back top",
    "strappy yoga bra",
    "sport cup",
    "strappy back bra",
    "strappy yoga top",
    "criss cross back",
    "long l",
    "low front bra",
    "gym bra",
    "workout sport bra",
    "cross back bralette",
    "pink criss cross top",
    "black strappy top",
    "cross line",
    "strappy workout bra",
    "back workout",
    "yoga bra top",
    "black strappy",
    "black criss cross top",
    "padded yoga bra",
    "pink bralette top",
    "green cross"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Reset the active selection to A1 (top-left), clearing the stale E15 selection
$null = $ws.Range("A1").Select()

